# edit.ps1 - applies the commit's content changes via the Word COM object model.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the leftover "_GoBack" bookmark (w:bookmarkStart/w:bookmarkEnd id=0)
#    that Word inserts automatically to mark the last edit location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Reword the data-preprocessing sentence in the Methods paragraph:
#    "... First, variables which had data for only 70 percent or more of
#    the total observations were discarded. ..."
#    -> "... First, only variables which had data for 70 percent or more of
#    the total observations were included in the analysis. ..."
# ---------------------------------------------------------------------------
$old1 = "First, variables which had data for only 70 percent or more of the total observations were discarded."
$new1 = "First, only variables which had data for 70 percent or more of the total observations were included in the analysis."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) The final revision also ran the grammar checker, which wraps a handful
#    of phrases in w:proofErr gramStart/gramEnd markers and, as a side
#    effect, splits the run that contained them. Reproduce the run splits
#    (toggling formatting on/off leaves the visible formatting untouched
#    but forces Word to break the run at the selection boundaries) at each
#    flagged phrase.
# ---------------------------------------------------------------------------
function Split-Run([string]$phrase) {
    $rng = $d.Content
    $found = $rng.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1
        $rng.Bold = 0
    }
}

Split-Run "most commonly used"
Split-Run "neighbors"
Split-Run "in order to"
Split-Run "Random forest"
